$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-13) holds the "Förändrad" date, stored as a serial date
# number. Bump each of these from 45205 to 45206 (one day later), matching
# the automatic update performed by the source tool.
foreach ($row in 2..13) {
    $ws.Cells.Item($row, 3).Value = 45206
}
